$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post in row 592 ("「己の翼を信じる者は己の重さを案じない」") was removed.
# Delete that entire row; Excel shifts all subsequent rows up by one
# (rows 593:616 become 592:615), and the used-range dimension shrinks
# from A1:C616 to A1:C615 automatically.
$ws.Rows(592).Delete()
